$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(3)

# --- New header cells L1:N1 ---
$ws.Range('L1').Value = 'email'
$ws.Range('M1').Value = 'lastName'
$ws.Range('N1').Value = 'firstName'

# --- Row 12 ---
$ws.Range('A12').Value = 'Mon Mar 07 2022'
$ws.Range('B12').Value = '07:54:48 GMT+0000 (Greenwich Mean Time)'
$ws.Range('D12').Value = 'User'
$ws.Range('E12').Value = '/api/auth/verify-email'
$ws.Range('F12').Value = 'login'
$ws.Range('G12').Value = 'failed'
$ws.Range('H12').Value = 'k0d3.s0n1k@gmail.com  login'
$ws.Range('I12').Value = 'error.invalid'
$ws.Range('L12').Value = 'k0d3.s0n1k@gmail.com'

# --- Row 13 ---
$ws.Range('A13').Value = 'Mon Mar 07 2022'
$ws.Range('B13').Value = '07:56:13 GMT+0000 (Greenwich Mean Time)'
$ws.Range('D13').Value = 'User'
$ws.Range('E13').Value = '/api/auth/verify-email'
$ws.Range('F13').Value = 'login'
$ws.Range('G13').Value = 'failed'
$ws.Range('H13').Value = 'k0d3.s0n1k@gmail.com  login'
$ws.Range('I13').Value = 'error.invalid'
$ws.Range('L13').Value = 'k0d3.s0n1k@gmail.com'

# --- Row 14 ---
$ws.Range('A14').Value = 'Mon Mar 07 2022'
$ws.Range('B14').Value = '07:57:30 GMT+0000 (Greenwich Mean Time)'
$ws.Range('D14').Value = 'User'
$ws.Range('E14').Value = '/api/auth/verify-email'
$ws.Range('F14').Value = 'login'
$ws.Range('G14').Value = 'failed'
$ws.Range('H14').Value = 'k0d3.s0n1k@gmail.com  login'
$ws.Range('I14').Value = 'error.invalid'
$ws.Range('L14').Value = 'k0d3.s0n1k@gmail.com'

# --- Row 15 ---
$ws.Range('A15').Value = 'Mon Mar 07 2022'
$ws.Range('B15').Value = '07:58:06 GMT+0000 (Greenwich Mean Time)'
$ws.Range('D15').Value = 'User'
$ws.Range('E15').Value = '/api/auth/login'
$ws.Range('F15').Value = 'login'
$ws.Range('G15').Value = 'failed'
$ws.Range('H15').Value = 'k0d3.s0n1k@gmail.com  login'
$ws.Range('I15').Value = 'error.invalid'
$ws.Range('L15').Value = 'k0d3.s0n1k@gmail.com'

# --- Row 16 ---
$ws.Range('A16').Value = 'Mon Mar 07 2022'
$ws.Range('B16').Value = '08:01:12 GMT+0000 (Greenwich Mean Time)'
$ws.Range('C16').Value = '''22892942601'
$ws.Range('D16').Value = 'User'
$ws.Range('E16').Value = '/api/auth/verify-otp'
$ws.Range('F16').Value = 'request'
$ws.Range('G16').Value = 'failed'
$ws.Range('H16').Value = '22892942601 request to receive otp'
$ws.Range('I16').Value = 'error.missing'

# --- Row 17 ---
$ws.Range('A17').Value = 'Mon Mar 07 2022'
$ws.Range('B17').Value = '08:02:08 GMT+0000 (Greenwich Mean Time)'
$ws.Range('C17').Value = '''22892942601'
$ws.Range('D17').Value = 'User'
$ws.Range('E17').Value = '/api/auth/verify-otp'
$ws.Range('F17').Value = 'request'
$ws.Range('G17').Value = 'succeeded'
$ws.Range('H17').Value = '22892942601 request to receive otp'
$ws.Range('J17').Value = 2
$ws.Range('K17').Value = '''2'

# --- Row 18 ---
$ws.Range('A18').Value = 'Mon Mar 07 2022'
$ws.Range('B18').Value = '08:05:54 GMT+0000 (Greenwich Mean Time)'
$ws.Range('D18').Value = 'User'
$ws.Range('E18').Value = '/api/auth/complete-infos'
$ws.Range('F18').Value = 'edit'
$ws.Range('G18').Value = 'succeeded'
$ws.Range('H18').Value = '    edit his infos'
$ws.Range('J18').Value = 2
$ws.Range('K18').Value = '''2'
$ws.Range('M18').Formula = '=""'
$ws.Range('N18').Formula = '=""'

# --- Row 19 ---
$ws.Range('A19').Value = 'Mon Mar 07 2022'
$ws.Range('B19').Value = '08:07:45 GMT+0000 (Greenwich Mean Time)'
$ws.Range('D19').Value = 'User'
$ws.Range('E19').Value = '/api/auth/profile'
$ws.Range('F19').Value = 'read'
$ws.Range('G19').Value = 'succeeded'
$ws.Range('H19').Value = 'Sonik  Kode  read his infos'
$ws.Range('J19').Value = 2
$ws.Range('K19').Value = '''2'
$ws.Range('M19').Value = 'Sonik'
$ws.Range('N19').Value = 'Kode'

# --- Row 20 ---
$ws.Range('A20').Value = 'Mon Mar 07 2022'
$ws.Range('B20').Value = '09:04:01 GMT+0000 (Greenwich Mean Time)'
$ws.Range('D20').Value = 'User'
$ws.Range('E20').Value = '/api/auth/add-email-auth'
$ws.Range('F20').Value = 'edit'
$ws.Range('G20').Value = 'failed'
$ws.Range('H20').Value = 'Sonik  Kode  edit his infos'
$ws.Range('I20').Value = '
Invalid `prisma.user.findFirst()` invocation:
{
  where: {
    id: 2,
    deletedAt: null
  },
  select: {
?   email?: true,
?   password?: true,
?   firstName?: true,
    lang: true,
    ~~~~
?   id?: true,
?   countryId?: true,
?   avatar?: true,
?   phoneNumber?: true,
?   emailVerifiedAt?: true,
?   phoneNumberVerifiedAt?: true,
?   lastName?: true,
?   birthDay?: true,
?   status?: true,
?   role?: true,
?   language?: true,
?   idCard?: true,
?   driverLicence?: true,
?   rating?: true,
?   createdAt?: true,
?   blockedAt?: true,
?   updatedAt?: true,
?   profileCompletedAt?: true,
?   deletedAt?: true,
?   deletionReport?: true,
?   operations?: true,
?   preferences?: true,
?   trips?: true,
?   travels?: true,
?   vehicles?: true,
?   historics?: true,
?   sendedNotifications?: true,
?   receivedNotifications?: true,
?   wallets?: true,
?   devices?: true,
?   country?: true,
?   _count?: true
  }
}
Unknown field `lang` for select statement on model User. Available options are listed in green. Did you mean `rating`?
'
$ws.Range('J20').Value = 2
$ws.Range('K20').Value = '''2'
$ws.Range('M20').Value = 'Sonik'
$ws.Range('N20').Value = 'Kode'

# --- Row 21 ---
$ws.Range('A21').Value = 'Mon Mar 07 2022'
$ws.Range('B21').Value = '09:09:07 GMT+0000 (Greenwich Mean Time)'
$ws.Range('D21').Value = 'User'
$ws.Range('E21').Value = '/api/auth/add-email-auth'
$ws.Range('F21').Value = 'edit'
$ws.Range('G21').Value = 'succeeded'
$ws.Range('H21').Value = 'Sonik  Kode  edit his infos'
$ws.Range('J21').Value = 2
$ws.Range('K21').Value = '''2'
$ws.Range('M21').Value = 'Sonik'
$ws.Range('N21').Value = 'Kode'

# --- Row 22 ---
$ws.Range('A22').Value = 'Mon Mar 07 2022'
$ws.Range('B22').Value = '09:12:57 GMT+0000 (Greenwich Mean Time)'
$ws.Range('D22').Value = 'User'
$ws.Range('E22').Value = '/api/auth/verify-email'
$ws.Range('F22').Value = 'login'
$ws.Range('G22').Value = 'failed'
$ws.Range('H22').Value = 'k0d3.s0n1k@gmail.com  login'
$ws.Range('I22').Value = 'error.unauthorized'
$ws.Range('L22').Value = 'k0d3.s0n1k@gmail.com'

# --- Row 23 ---
$ws.Range('A23').Value = 'Mon Mar 07 2022'
$ws.Range('B23').Value = '09:28:45 GMT+0000 (Greenwich Mean Time)'
$ws.Range('D23').Value = 'User'
$ws.Range('E23').Value = '/api/auth/confirm-email'
$ws.Range('F23').Value = 'activate'
$ws.Range('G23').Value = 'succeeded'
$ws.Range('H23').Value = 'k0d3.s0n1k@gmail.com activate his account'
$ws.Range('K23').Formula = '=""'
$ws.Range('L23').Value = 'k0d3.s0n1k@gmail.com'

# --- Row 24 ---
$ws.Range('A24').Value = 'Mon Mar 07 2022'
$ws.Range('B24').Value = '09:34:56 GMT+0000 (Greenwich Mean Time)'
$ws.Range('D24').Value = 'User'
$ws.Range('E24').Value = '/api/auth/verify-email'
$ws.Range('F24').Value = 'login'
$ws.Range('G24').Value = 'succeeded'
$ws.Range('H24').Value = 'k0d3.s0n1k@gmail.com  login'
$ws.Range('L24').Value = 'k0d3.s0n1k@gmail.com'

# --- Row 25 ---
$ws.Range('A25').Value = 'Mon Mar 07 2022'
$ws.Range('B25').Value = '09:35:47 GMT+0000 (Greenwich Mean Time)'
$ws.Range('D25').Value = 'User'
$ws.Range('E25').Value = '/api/auth/verify-otp'
$ws.Range('F25').Value = 'request'
$ws.Range('G25').Value = 'failed'
$ws.Range('H25').Value = 'undefined request to receive otp'
$ws.Range('I25').Value = 'error.missing'

# --- Row 26 ---
$ws.Range('A26').Value = 'Mon Mar 07 2022'
$ws.Range('B26').Value = '09:37:54 GMT+0000 (Greenwich Mean Time)'
$ws.Range('D26').Value = 'User'
$ws.Range('E26').Value = '/api/auth/verify-password'
$ws.Range('F26').Value = 'login'
$ws.Range('G26').Value = 'failed'
$ws.Range('H26').Value = '2  login'
$ws.Range('I26').Value = 'error.invalid'
$ws.Range('J26').Value = 2
$ws.Range('K26').Value = 2

# --- Row 27 ---
$ws.Range('A27').Value = 'Mon Mar 07 2022'
$ws.Range('B27').Value = '09:38:05 GMT+0000 (Greenwich Mean Time)'
$ws.Range('D27').Value = 'User'
$ws.Range('E27').Value = '/api/auth/verify-email'
$ws.Range('F27').Value = 'login'
$ws.Range('G27').Value = 'succeeded'
$ws.Range('H27').Value = 'k0d3.s0n1k@gmail.com  login'
$ws.Range('L27').Value = 'k0d3.s0n1k@gmail.com'

# --- Row 28 ---
$ws.Range('A28').Value = 'Mon Mar 07 2022'
$ws.Range('B28').Value = '09:38:31 GMT+0000 (Greenwich Mean Time)'
$ws.Range('D28').Value = 'User'
$ws.Range('E28').Value = '/api/auth/verify-password'
$ws.Range('F28').Value = 'login'
$ws.Range('G28').Value = 'failed'
$ws.Range('H28').Value = '2  login'
$ws.Range('I28').Value = 'error.invalid'
$ws.Range('J28').Value = 2
$ws.Range('K28').Value = 2

# --- Row 29 ---
$ws.Range('A29').Value = 'Mon Mar 07 2022'
$ws.Range('B29').Value = '09:38:52 GMT+0000 (Greenwich Mean Time)'
$ws.Range('D29').Value = 'User'
$ws.Range('E29').Value = '/api/auth/verify-password'
$ws.Range('F29').Value = 'login'
$ws.Range('G29').Value = 'failed'
$ws.Range('H29').Value = '2  login'
$ws.Range('I29').Value = 'error.invalid'
$ws.Range('J29').Value = 2
$ws.Range('K29').Value = 2

# --- Row 30 ---
$ws.Range('A30').Value = 'Mon Mar 07 2022'
$ws.Range('B30').Value = '09:40:34 GMT+0000 (Greenwich Mean Time)'
$ws.Range('D30').Value = 'User'
$ws.Range('E30').Value = '/api/auth/verify-email'
$ws.Range('F30').Value = 'login'
$ws.Range('G30').Value = 'succeeded'
$ws.Range('H30').Value = 'k0d3.s0n1k@gmail.com  login'
$ws.Range('L30').Value = 'k0d3.s0n1k@gmail.com'

# --- Row 31 ---
$ws.Range('A31').Value = 'Mon Mar 07 2022'
$ws.Range('B31').Value = '09:43:02 GMT+0000 (Greenwich Mean Time)'
$ws.Range('C31').Value = '''22892942601'
$ws.Range('D31').Value = 'User'
$ws.Range('E31').Value = '/api/auth/send-otp'
$ws.Range('F31').Value = 'request'
$ws.Range('G31').Value = 'succeeded'
$ws.Range('H31').Value = '22892942601 request to receive otp'

# --- Row 32 ---
$ws.Range('A32').Value = 'Mon Mar 07 2022'
$ws.Range('B32').Value = '09:45:09 GMT+0000 (Greenwich Mean Time)'
$ws.Range('D32').Value = 'User'
$ws.Range('E32').Value = '/api/auth/verify-password'
$ws.Range('F32').Value = 'login'
$ws.Range('G32').Value = 'failed'
$ws.Range('H32').Value = '2  login'
$ws.Range('I32').Value = 'error.invalid'
$ws.Range('J32').Value = 2
$ws.Range('K32').Value = 2

# --- Row 33 ---
$ws.Range('A33').Value = 'Mon Mar 07 2022'
$ws.Range('B33').Value = '09:45:43 GMT+0000 (Greenwich Mean Time)'
$ws.Range('C33').Value = '''22892942601'
$ws.Range('D33').Value = 'User'
$ws.Range('E33').Value = '/api/auth/verify-otp'
$ws.Range('F33').Value = 'request'
$ws.Range('G33').Value = 'failed'
$ws.Range('H33').Value = '22892942601 request to receive otp'
$ws.Range('I33').Value = 'error.invalid'

# --- Row 34 ---
$ws.Range('A34').Value = 'Mon Mar 07 2022'
$ws.Range('B34').Value = '09:48:44 GMT+0000 (Greenwich Mean Time)'
$ws.Range('D34').Value = 'User'
$ws.Range('E34').Value = '/api/auth/verify-password'
$ws.Range('F34').Value = 'login'
$ws.Range('G34').Value = 'failed'
$ws.Range('H34').Value = '2  login'
$ws.Range('I34').Value = 'error.invalid'
$ws.Range('J34').Value = 2
$ws.Range('K34').Value = 2

# --- Row 35 ---
$ws.Range('A35').Value = 'Mon Mar 07 2022'
$ws.Range('B35').Value = '09:48:51 GMT+0000 (Greenwich Mean Time)'
$ws.Range('D35').Value = 'User'
$ws.Range('E35').Value = '/api/auth/verify-email'
$ws.Range('F35').Value = 'login'
$ws.Range('G35').Value = 'succeeded'
$ws.Range('H35').Value = 'k0d3.s0n1k@gmail.com  login'
$ws.Range('L35').Value = 'k0d3.s0n1k@gmail.com'

# --- Row 36 ---
$ws.Range('A36').Value = 'Mon Mar 07 2022'
$ws.Range('B36').Value = '09:49:17 GMT+0000 (Greenwich Mean Time)'
$ws.Range('D36').Value = 'User'
$ws.Range('E36').Value = '/api/auth/verify-password'
$ws.Range('F36').Value = 'login'
$ws.Range('G36').Value = 'failed'
$ws.Range('H36').Value = '2  login'
$ws.Range('I36').Value = 'error.invalid'
$ws.Range('J36').Value = 2
$ws.Range('K36').Value = 2

# --- Row 37 ---
$ws.Range('A37').Value = 'Mon Mar 07 2022'
$ws.Range('B37').Value = '19:19:47 GMT+0000 (Greenwich Mean Time)'
$ws.Range('D37').Value = 'User'
$ws.Range('E37').Value = '/api/auth/verify-email'
$ws.Range('F37').Value = 'login'
$ws.Range('G37').Value = 'failed'
$ws.Range('H37').Value = 'k0d3.s0n1k@gmail.com  login'
$ws.Range('I37').Value = '
Invalid `prisma.user.findFirst()` invocation:
  Can''t reach database server at `ec2-54-216-17-9.eu-west-1.compute.amazonaws.com`:`5432`
Please make sure your database server is running at `ec2-54-216-17-9.eu-west-1.compute.amazonaws.com`:`5432`.'
$ws.Range('L37').Value = 'k0d3.s0n1k@gmail.com'

# --- Row 38 ---
$ws.Range('A38').Value = 'Mon Mar 07 2022'
$ws.Range('B38').Value = '19:22:31 GMT+0000 (Greenwich Mean Time)'
$ws.Range('D38').Value = 'User'
$ws.Range('E38').Value = '/api/auth/verify-password'
$ws.Range('F38').Value = 'login'
$ws.Range('G38').Value = 'failed'
$ws.Range('H38').Value = '2  login'
$ws.Range('I38').Value = 'error.invalid'
$ws.Range('J38').Value = 2
$ws.Range('K38').Value = 2

# --- Row 39 ---
$ws.Range('A39').Value = 'Mon Mar 07 2022'
$ws.Range('B39').Value = '19:22:36 GMT+0000 (Greenwich Mean Time)'
$ws.Range('D39').Value = 'User'
$ws.Range('E39').Value = '/api/auth/verify-email'
$ws.Range('F39').Value = 'login'
$ws.Range('G39').Value = 'failed'
$ws.Range('H39').Value = 'k0d3.s0n1k@gmail.com  login'
$ws.Range('I39').Value = '
Invalid `prisma.user.findFirst()` invocation:
  Can''t reach database server at `ec2-54-216-17-9.eu-west-1.compute.amazonaws.com`:`5432`
Please make sure your database server is running at `ec2-54-216-17-9.eu-west-1.compute.amazonaws.com`:`5432`.'
$ws.Range('L39').Value = 'k0d3.s0n1k@gmail.com'

# --- Row 40 ---
$ws.Range('A40').Value = 'Mon Mar 07 2022'
$ws.Range('B40').Value = '22:51:37 GMT+0000 (Greenwich Mean Time)'
$ws.Range('D40').Value = 'User'
$ws.Range('E40').Value = '/api/auth/verify-email'
$ws.Range('F40').Value = 'login'
$ws.Range('G40').Value = 'failed'
$ws.Range('H40').Value = 'k0d3.s0n1k@gmail.com  login'
$ws.Range('I40').Value = '
Invalid `prisma.user.findFirst()` invocation:
  Can''t reach database server at `ec2-54-216-17-9.eu-west-1.compute.amazonaws.com`:`5432`
Please make sure your database server is running at `ec2-54-216-17-9.eu-west-1.compute.amazonaws.com`:`5432`.'
$ws.Range('L40').Value = 'k0d3.s0n1k@gmail.com'

# --- Row 41 ---
$ws.Range('A41').Value = 'Mon Mar 07 2022'
$ws.Range('B41').Value = '22:54:10 GMT+0000 (Greenwich Mean Time)'
$ws.Range('D41').Value = 'User'
$ws.Range('E41').Value = '/api/auth/verify-email'
$ws.Range('F41').Value = 'login'
$ws.Range('G41').Value = 'succeeded'
$ws.Range('H41').Value = 'k0d3.s0n1k@gmail.com  login'
$ws.Range('L41').Value = 'k0d3.s0n1k@gmail.com'

# --- Row 42 ---
$ws.Range('A42').Value = 'Mon Mar 07 2022'
$ws.Range('B42').Value = '22:56:48 GMT+0000 (Greenwich Mean Time)'
$ws.Range('D42').Value = 'User'
$ws.Range('E42').Value = '/api/auth/verify-password'
$ws.Range('F42').Value = 'login'
$ws.Range('G42').Value = 'failed'
$ws.Range('H42').Value = '2  login'
$ws.Range('I42').Value = 'error.invalid'
$ws.Range('J42').Value = 2
$ws.Range('K42').Value = 2

# --- Row 43 ---
$ws.Range('A43').Value = 'Mon Mar 07 2022'
$ws.Range('B43').Value = '22:58:40 GMT+0000 (Greenwich Mean Time)'
$ws.Range('D43').Value = 'User'
$ws.Range('E43').Value = '/api/auth/verify-password'
$ws.Range('F43').Value = 'login'
$ws.Range('G43').Value = 'failed'
$ws.Range('H43').Value = '2  login'
$ws.Range('I43').Value = 'error.invalid'
$ws.Range('J43').Value = 2
$ws.Range('K43').Value = 2
